# New crime data collected - weekly CompStat report update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/20/2023  Through  11/26/2023"

# --- Row 14 (Murder) needs D/E switched from text placeholders to numbers ---
$ws.Range("D14").NumberFormat = $ws.Range("D15").NumberFormat
$ws.Range("E14").NumberFormat = $ws.Range("E15").NumberFormat

# --- Row 30 (Hate Crimes) needs C switched from text placeholder to number ---
$ws.Range("C30").NumberFormat = $ws.Range("C29").NumberFormat

# --- Week to Date / 28 Day / Year to Date / 2 Year table (rows 14-30) ---

# Row 14: Murder
$ws.Range("C14").Value = "0"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 4
$ws.Range("H14").Value = -25
$ws.Range("I14").Value = 53
$ws.Range("J14").Value = 59
$ws.Range("K14").Value = -10.169491525423
$ws.Range("L14").Value = 10.416666666666
$ws.Range("M14").Value = -31.168831168831
$ws.Range("N14").Value = -76.444444444444

# Row 15: Rape
$ws.Range("A15").Value = "Rape"
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 50
$ws.Range("F15").Value = 21
$ws.Range("G15").Value = 17
$ws.Range("H15").Value = 23.529411764705
$ws.Range("I15").Value = 200
$ws.Range("J15").Value = 196
$ws.Range("K15").Value = 2.040816326530
$ws.Range("L15").Value = 2.040816326530
$ws.Range("M15").Value = 21.212121212121
$ws.Range("N15").Value = -60.784313725490

# Row 16: Robbery
$ws.Range("A16").Value = "Robbery"
$ws.Range("C16").Value = 33
$ws.Range("D16").Value = 36
$ws.Range("E16").Value = -8.333333333333
$ws.Range("F16").Value = 148
$ws.Range("G16").Value = 152
$ws.Range("H16").Value = -2.631578947368
$ws.Range("I16").Value = 1610
$ws.Range("J16").Value = 1775
$ws.Range("K16").Value = -9.295774647887
$ws.Range("L16").Value = 24.903025601241
$ws.Range("M16").Value = -38.992042440318
$ws.Range("N16").Value = -87.339781394983

# Row 17: Fel. Assault
$ws.Range("A17").Value = "Fel. Assault"
$ws.Range("C17").Value = 80
$ws.Range("D17").Value = 60
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 264
$ws.Range("G17").Value = 250
$ws.Range("H17").Value = 5.6
$ws.Range("I17").Value = 3200
$ws.Range("J17").Value = 3126
$ws.Range("K17").Value = 2.367242482405
$ws.Range("L17").Value = 12.715744980627
$ws.Range("M17").Value = 42.538975501113
$ws.Range("N17").Value = -47.780678851174

# Row 18: Burglary
$ws.Range("A18").Value = "Burglary"
$ws.Range("C18").Value = 24
$ws.Range("D18").Value = 40
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 100
$ws.Range("G18").Value = 171
$ws.Range("H18").Value = -41.520467836257
$ws.Range("I18").Value = 1526
$ws.Range("J18").Value = 1957
$ws.Range("K18").Value = -22.023505365355
$ws.Range("L18").Value = -4.803493449781
$ws.Range("M18").Value = -51.830808080808
$ws.Range("N18").Value = -90.936627665261

# Row 19: Gr. Larceny
$ws.Range("A19").Value = "Gr. Larceny"
$ws.Range("C19").Value = 84
$ws.Range("D19").Value = 115
$ws.Range("E19").Value = -26.956521739130
$ws.Range("F19").Value = 452
$ws.Range("G19").Value = 518
$ws.Range("H19").Value = -12.741312741312
$ws.Range("I19").Value = 5867
$ws.Range("J19").Value = 6509
$ws.Range("K19").Value = -9.863266246735
$ws.Range("L19").Value = 25.202731540759
$ws.Range("M19").Value = 17.528044871794
$ws.Range("N19").Value = -27.397599307016

# Row 20: G.L.A.
$ws.Range("A20").Value = "G.L.A."
$ws.Range("C20").Value = 38
$ws.Range("D20").Value = 27
$ws.Range("E20").Value = 40.740740740740
$ws.Range("F20").Value = 150
$ws.Range("G20").Value = 113
$ws.Range("H20").Value = 32.743362831858
$ws.Range("I20").Value = 1699
$ws.Range("J20").Value = 1637
$ws.Range("K20").Value = 3.787416004886
$ws.Range("L20").Value = 43.013468013468
$ws.Range("M20").Value = -3.080433542498
$ws.Range("N20").Value = -91.957777146643

# Row 21: TOTAL
$ws.Range("A21").Value = "TOTAL"
$ws.Range("C21").Value = 265
$ws.Range("D21").Value = 283
$ws.Range("E21").Value = -6.360424028268
$ws.Range("F21").Value = 1138
$ws.Range("G21").Value = 1225
$ws.Range("H21").Value = -7.102040816326
$ws.Range("I21").Value = 14155
$ws.Range("J21").Value = 15259
$ws.Range("K21").Value = -7.235074382331
$ws.Range("L21").Value = 19.461557937378
$ws.Range("M21").Value = -5.878050402287
$ws.Range("N21").Value = -78.430147507009

# Row 22: Transit
$ws.Range("A22").Value = "Transit"
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 24
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = 71.428571428571
$ws.Range("I22").Value = 179
$ws.Range("J22").Value = 175
$ws.Range("K22").Value = 2.285714285714
$ws.Range("L22").Value = 27.857142857142
$ws.Range("M22").Value = -29.249011857707
$ws.Range("N22").Value = "***.*"

# Row 23: Housing
$ws.Range("C23").Value = 10
$ws.Range("D23").Value = 12
$ws.Range("E23").Value = -16.666666666666
$ws.Range("F23").Value = 33
$ws.Range("G23").Value = 48
$ws.Range("H23").Value = -31.25
$ws.Range("I23").Value = 473
$ws.Range("J23").Value = 512
$ws.Range("K23").Value = -7.6171875
$ws.Range("L23").Value = 8.486238532110
$ws.Range("M23").Value = 52.090032154340
$ws.Range("N23").Value = "***.*"

# Row 24: Petit Larceny
$ws.Range("C24").Value = 227
$ws.Range("D24").Value = 276
$ws.Range("E24").Value = -17.753623188405
$ws.Range("F24").Value = 1140
$ws.Range("G24").Value = 1276
$ws.Range("H24").Value = -10.658307210031
$ws.Range("I24").Value = 14343
$ws.Range("J24").Value = 14790
$ws.Range("K24").Value = -3.022312373225
$ws.Range("L24").Value = 31.575084854600
$ws.Range("M24").Value = 27.075396473819
$ws.Range("N24").Value = "***.*"

# Row 25: Misd. Assault
$ws.Range("C25").Value = 112
$ws.Range("D25").Value = 116
$ws.Range("E25").Value = -3.448275862068
$ws.Range("F25").Value = 469
$ws.Range("G25").Value = 425
$ws.Range("H25").Value = 10.352941176470
$ws.Range("I25").Value = 5348
$ws.Range("J25").Value = 5049
$ws.Range("K25").Value = 5.921964745494
$ws.Range("L25").Value = 18.686196182867
$ws.Range("M25").Value = -12.714215766280
$ws.Range("N25").Value = "***.*"

# Row 26: UCR Rape*
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 34
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 36
$ws.Range("I26").Value = 294
$ws.Range("J26").Value = 315
$ws.Range("K26").Value = -6.666666666666
$ws.Range("L26").Value = -5.769230769230
$ws.Range("M26").Value = "***.*"
$ws.Range("N26").Value = "***.*"

# Row 27: Other Sex Crimes
$ws.Range("C27").Value = 14
$ws.Range("D27").Value = 9
$ws.Range("E27").Value = 55.555555555555
$ws.Range("F27").Value = 51
$ws.Range("G27").Value = 53
$ws.Range("H27").Value = -3.773584905660
$ws.Range("I27").Value = 585
$ws.Range("J27").Value = 636
$ws.Range("K27").Value = -8.018867924528
$ws.Range("L27").Value = 0.171232876712
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"

# Row 28: Shooting Vic.
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 300
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 11.111111111111
$ws.Range("I28").Value = 134
$ws.Range("J28").Value = 197
$ws.Range("K28").Value = -31.979695431472
$ws.Range("L28").Value = -27.567567567567
$ws.Range("M28").Value = -47.65625
$ws.Range("N28").Value = -81.232492997198

# Row 29: Shooting Inc.
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 100
$ws.Range("F29").Value = 8
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 115
$ws.Range("J29").Value = 155
$ws.Range("K29").Value = -25.806451612903
$ws.Range("L29").Value = -29.447852760736
$ws.Range("M29").Value = -45.497630331753
$ws.Range("N29").Value = -81.629392971246

# Row 30: Hate Crimes
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 15
$ws.Range("G30").Value = 11
$ws.Range("H30").Value = 36.363636363636
$ws.Range("I30").Value = 100
$ws.Range("J30").Value = 112
$ws.Range("K30").Value = -10.714285714285
$ws.Range("L30").Value = 47.058823529411
$ws.Range("M30").Value = "***.*"
$ws.Range("N30").Value = "***.*"
